$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "sssss"
$ws.Range("C4").Value = "ffffff"
$ws.Range("E4").Value = 234678902

$ws.Range("E4").Select()
